# Automatische test-sync: 2025-08-14 21:26:50
# Adds one new mail-log entry ("Demo inplannen" / INTERN - Planning / Afspraak)
# to the Logs sheet, and updates the Dashboard category summary + chart
# ranges accordingly.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# --- Logs sheet: append the new row 26 -------------------------------------
$logs.Cells.Item(26, 1).Value = "Demo inplannen"
$logs.Cells.Item(26, 2).Value = "klantenservice@testbedrijf123.nl"
$logs.Cells.Item(26, 3).Value = "Kun je vrijdag om 11:00 een demo inplannen bij Van Dijk?"
$logs.Cells.Item(26, 4).Value = "INTERN – Planning / Afspraak"
$logs.Cells.Item(26, 5).Value = "Bedankt, we hebben dit doorgestuurd naar planning@testbedrijf123.nl."
$logs.Cells.Item(26, 6).Value = "2025-08-14 21:25:51"
$logs.Cells.Item(26, 7).Value = "Nee"
$logs.Cells.Item(26, 8).Value = "Ja"
$logs.Cells.Item(26, 9).Value = "Nee"
$logs.Cells.Item(26, 10).Value = "Nee"

# --- Dashboard sheet: append the new category row 8 -------------------------
$dash.Cells.Item(8, 1).Value = "INTERN – Planning / Afspraak"
$dash.Cells.Item(8, 2).Value = 1

# --- Logs sheet: extend conditional formatting ranges to include row 26 -----
$newD = $logs.Range("D2:D26")
$fcsD = $logs.Range("D2:D25").FormatConditions
for ($i = 1; $i -le $fcsD.Count; $i++) {
    $fcsD.Item($i).ModifyAppliesToRange($newD)
}

$newG = $logs.Range("G2:G26")
$fcsG = $logs.Range("G2:G25").FormatConditions
for ($i = 1; $i -le $fcsG.Count; $i++) {
    $fcsG.Item($i).ModifyAppliesToRange($newG)
}

$newH = $logs.Range("H2:H26")
$fcsH = $logs.Range("H2:H25").FormatConditions
for ($i = 1; $i -le $fcsH.Count; $i++) {
    $fcsH.Item($i).ModifyAppliesToRange($newH)
}

$newI = $logs.Range("I2:I26")
$fcsI = $logs.Range("I2:I25").FormatConditions
for ($i = 1; $i -le $fcsI.Count; $i++) {
    $fcsI.Item($i).ModifyAppliesToRange($newI)
}

$newJ = $logs.Range("J2:J26")
$fcsJ = $logs.Range("J2:J25").FormatConditions
for ($i = 1; $i -le $fcsJ.Count; $i++) {
    $fcsJ.Item($i).ModifyAppliesToRange($newJ)
}

# --- Dashboard chart: extend the category/value series ranges to row 8 ------
$chartObj = $dash.ChartObjects().Item(1)
$chart = $chartObj.Chart
$ser = $chart.SeriesCollection().Item(1)
$ser.Formula = "=SERIES(Dashboard!B1,Dashboard!`$A`$2:`$A`$8,Dashboard!`$B`$2:`$B`$8,1)"

